# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 07:16"

# 2. Swap the "Almeria" / "Lugo" rows (row 47 <-> row 48):
#    - Province names in column A swap
#    - "Casos activos" values in column C swap
#    (Casos totales / Recuperados / Muertes are identical between the two
#    rows, so no visible change is needed there.)
$nameA47 = $ws.Range("A47").Value2
$nameA48 = $ws.Range("A48").Value2
$ws.Range("A47").Value2 = $nameA48
$ws.Range("A48").Value2 = $nameA47

$valC47 = $ws.Range("C47").Value2
$valC48 = $ws.Range("C48").Value2
$ws.Range("C47").Value2 = $valC48
$ws.Range("C48").Value2 = $valC47
